$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 966.6
$ws.Range("J2").Value = 3000.6667
$ws.Range("L2").Value = 3000.6667
$ws.Range("N2").Value = -3226.6667
$ws.Range("H6").Value = 17755.87
$ws.Range("I6").Value = 22543.723
$ws.Range("J6").Value = 519.6
$ws.Range("K6").Value = 67631.16900000001
$ws.Range("L6").Value = 1558.8
$ws.Range("M6").Value = -67519.16900000001
$ws.Range("N6").Value = -1782.8
$ws.Range("H17").Value = 2040.3556
$ws.Range("J17").Value = 2040.3556
$ws.Range("L17").Value = 6121.066800000001
$ws.Range("N17").Value = -6457.066800000001
$ws.Range("H28").Value = 1627.9286
$ws.Range("I28").Value = 1456.7778
$ws.Range("K28").Value = 1456.7778
$ws.Range("M28").Value = -971.7778000000001
$ws.Range("H51").Value = 9457.385
$ws.Range("I51").Value = 7183
$ws.Range("J51").Value = 10139.7
$ws.Range("K51").Value = 7183
$ws.Range("L51").Value = 10139.7
$ws.Range("M51").Value = -6699
$ws.Range("N51").Value = -11107.7
$ws.Range("H58").Value = 780.7143
$ws.Range("I58").Value = 79.833336
$ws.Range("J58").Value = 4986
$ws.Range("K58").Value = 239.500008
$ws.Range("L58").Value = 14958
$ws.Range("M58").Value = -89.50000800000001
$ws.Range("N58").Value = -15258
$ws.Range("H103").Value = 872.0909
$ws.Range("J103").Value = 932.55554
$ws.Range("L103").Value = 2797.66662
$ws.Range("N103").Value = -3969.66662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 758.1111
$ws.Range("I2").Value = 683
$ws.Range("J2").Value = 972.7143
$ws.Range("K2").Value = 683
$ws.Range("L2").Value = 972.7143
$ws.Range("M2").Value = -570
$ws.Range("N2").Value = -1198.7143
$ws.Range("H32").Value = 4423.337
$ws.Range("I32").Value = 2046.4756
$ws.Range("K32").Value = 2046.4756
$ws.Range("M32").Value = -1759.4756
$ws.Range("H63").Value = 3633.4285
$ws.Range("I63").Value = 2119.75
$ws.Range("K63").Value = 2119.75
$ws.Range("M63").Value = -1433.75
$ws.Range("H66").Value = 3633.4285
$ws.Range("I66").Value = 2119.75
$ws.Range("K66").Value = 10598.75
$ws.Range("M66").Value = -7166.75
$ws.Range("H88").Value = 2158.1667
$ws.Range("I88").Value = 1646.6
$ws.Range("J88").Value = 2523.5715
$ws.Range("K88").Value = 1646.6
$ws.Range("L88").Value = 2523.5715
$ws.Range("M88").Value = -1240.6
$ws.Range("N88").Value = -3335.5715
$ws.Range("H91").Value = 2158.1667
$ws.Range("I91").Value = 1646.6
$ws.Range("J91").Value = 2523.5715
$ws.Range("K91").Value = 1646.6
$ws.Range("L91").Value = 2523.5715
$ws.Range("M91").Value = -242.5999999999999
$ws.Range("N91").Value = -5331.5715
$ws.Range("H116").Value = 758.1111
$ws.Range("I116").Value = 683
$ws.Range("J116").Value = 972.7143
$ws.Range("K116").Value = 683
$ws.Range("L116").Value = 972.7143
$ws.Range("M116").Value = 1611
$ws.Range("N116").Value = -5560.7143
$ws.Range("H122").Value = 3234.1667
$ws.Range("I122").Value = 2061.04
$ws.Range("K122").Value = 6183.12
$ws.Range("M122").Value = -3733.12

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 758.1111
$ws.Range("I3").Value = 683
$ws.Range("J3").Value = 972.7143
$ws.Range("K3").Value = 683
$ws.Range("L3").Value = 972.7143
$ws.Range("M3").Value = -569
$ws.Range("N3").Value = -1200.7143
$ws.Range("H46").Value = 52500
$ws.Range("J46").Value = 52500
$ws.Range("L46").Value = 52500
$ws.Range("N46").Value = -53096
$ws.Range("H94").Value = 1128.5
$ws.Range("I94").Value = 841.75
$ws.Range("K94").Value = 841.75
$ws.Range("M94").Value = -390.75
$ws.Range("H128").Value = 100000
$ws.Range("I128").Value = 100000
$ws.Range("K128").Value = 300000
$ws.Range("M128").Value = -297510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3555.5557
$ws.Range("I122").Value = 3142.8572
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9428.571599999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6978.571599999999
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 4785.5
$ws.Range("I132").Value = 2993.2666
$ws.Range("J132").Value = 10162.2
$ws.Range("K132").Value = 8979.799800000001
$ws.Range("L132").Value = 30486.6
$ws.Range("M132").Value = -6449.799800000001
$ws.Range("N132").Value = -35546.60000000001
$ws.Range("H134").Value = 5228.636
$ws.Range("I134").Value = 4187.625
$ws.Range("J134").Value = 8004.6665
$ws.Range("K134").Value = 12562.875
$ws.Range("L134").Value = 24013.9995
$ws.Range("M134").Value = -10027.875
$ws.Range("N134").Value = -29083.9995
$ws.Range("H138").Value = 80418.42999999999
$ws.Range("J138").Value = 80418.42999999999
$ws.Range("L138").Value = 80418.42999999999
$ws.Range("N138").Value = -90698.42999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1050
$ws.Range("J113").Value = 1165.6666
$ws.Range("L113").Value = 3496.9998
$ws.Range("N113").Value = -7836.9998
$ws.Range("H128").Value = 459998.4
$ws.Range("I128").Value = 459998.4
$ws.Range("K128").Value = 1379995.2
$ws.Range("M128").Value = -1375015.2
$ws.Range("H131").Value = 829262.9399999999
$ws.Range("I131").Value = 970.7273
$ws.Range("J131").Value = 1657555.1
$ws.Range("K131").Value = 2912.1819
$ws.Range("L131").Value = 4972665.300000001
$ws.Range("M131").Value = 2127.8181
$ws.Range("N131").Value = -4982745.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 29999
$ws.Range("J38").Value = 29999
$ws.Range("L38").Value = 29999
$ws.Range("N38").Value = -30925
$ws.Range("H46").Value = 40000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 40000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 40000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -40312
$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 5000
$ws.Range("K57").Value = 5000
$ws.Range("M57").Value = -4180
$ws.Range("H58").Value = 29500
$ws.Range("I58").Value = 40000
$ws.Range("J58").Value = 19000
$ws.Range("K58").Value = 40000
$ws.Range("L58").Value = 19000
$ws.Range("M58").Value = -39723
$ws.Range("N58").Value = -19554
$ws.Range("H80").Value = 12985.818
$ws.Range("I80").Value = 13022.615
$ws.Range("J80").Value = 12932.667
$ws.Range("K80").Value = 13022.615
$ws.Range("L80").Value = 12932.667
$ws.Range("M80").Value = -12024.615
$ws.Range("N80").Value = -14928.667
$ws.Range("H83").Value = 12985.818
$ws.Range("I83").Value = 13022.615
$ws.Range("J83").Value = 12932.667
$ws.Range("K83").Value = 65113.075
$ws.Range("L83").Value = 64663.335
$ws.Range("M83").Value = -60121.075
$ws.Range("N83").Value = -74647.33499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8022.6
$ws.Range("I40").Value = 6737.625
$ws.Range("J40").Value = 13162.5
$ws.Range("K40").Value = 6737.625
$ws.Range("L40").Value = 13162.5
$ws.Range("M40").Value = -6601.625
$ws.Range("N40").Value = -13434.5
$ws.Range("H100").Value = 5750
$ws.Range("I100").Value = 8500
$ws.Range("K100").Value = 8500
$ws.Range("M100").Value = -7959
$ws.Range("H122").Value = 6114.421
$ws.Range("I122").Value = 4845
$ws.Range("K122").Value = 14535
$ws.Range("M122").Value = -12085
$ws.Range("H132").Value = 3281
$ws.Range("I132").Value = 2355.4666
$ws.Range("J132").Value = 6057.6
$ws.Range("K132").Value = 7066.399800000001
$ws.Range("L132").Value = 18172.8
$ws.Range("M132").Value = -4536.399800000001
$ws.Range("N132").Value = -23232.8
$ws.Range("H134").Value = 50166.668
$ws.Range("J134").Value = 50166.668
$ws.Range("L134").Value = 50166.668
$ws.Range("N134").Value = -60306.668
$ws.Range("H136").Value = 4522.25
$ws.Range("I136").Value = 2293.9375
$ws.Range("J136").Value = 6304.9
$ws.Range("K136").Value = 6881.8125
$ws.Range("L136").Value = 18914.7
$ws.Range("M136").Value = -4331.8125
$ws.Range("N136").Value = -24014.7
$ws.Range("H140").Value = 58639.1
$ws.Range("J140").Value = 58639.1
$ws.Range("L140").Value = 58639.1
$ws.Range("N140").Value = -68999.10000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 748.6667
$ws.Range("I100").Value = 493
$ws.Range("J100").Value = 1004.3333
$ws.Range("K100").Value = 986
$ws.Range("L100").Value = 2008.6666
$ws.Range("M100").Value = -445
$ws.Range("N100").Value = -3090.6666
